$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signup_CS")

# Update trans_ref for the first test case
$ws.Range("C3").Value = "20220706-01"

# Update the API response fields to reflect a QR match failure instead of success
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "2018"
$ws.Range("K3").Value = "QR match fail"
$ws.Range("M3").Value = "E0012"
$ws.Range("N3").Value = "The information of ID card and QR code is not matched."
$ws.Range("O3").Value = "ข้อมูลบัตรประจำตัวประชาชนและ QR Code ไม่ตรงกัน"
$ws.Range("P3").Value = '{"status":{"code":"2018","message":"QR match fail","remark":"","user_code":"E0012","user_message_en":"The information of ID card and QR code is not matched.","user_message_th":"ข้อมูลบัตรประจำตัวประชาชนและ QR Code ไม่ตรงกัน"},"data":null}'

# Move the active selection to D5 (matches the saved cursor position in the workbook)
$ws.Range("D5").Select()
